$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The last paragraph currently reads "TB_Strona_name-tab" and carries
# the (hidden) "_GoBack" bookmark between "Strona" and the second "_".
# The edit:
#   1) removes that bookmark from its old position, and
#   2) appends a brand new paragraph "HL_Strona_name-hyperlink" right
#      after it, with "Strona" and "name" italicised (matching the
#      style used throughout the document) and with the "_GoBack"
#      bookmark now sitting at the very end of that new paragraph.
# ------------------------------------------------------------------

# Step 1: relocate the "_GoBack" bookmark by deleting it from its
# current spot (it will be recreated at the end of the new paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Step 2: find the last paragraph ("...-tab") and open a new, empty
# paragraph right after it, inheriting its paragraph formatting.
$tabParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$tabRange = $tabParagraph.Range
$tabRange.Collapse(0)
$tabRange.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newParagraph.Range

# Step 3: populate the new paragraph via an OOXML fragment so that each
# differently-formatted piece of text ("HL", "_", "Strona", "_", "name",
# "-hyperlink") stays in its own run, exactly like the rest of the
# document's "<ABC>_Strona_name-..." lines, and append the "_GoBack"
# bookmark at the end.
$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr>
<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>HL</w:t></w:r>
<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>_</w:t></w:r>
<w:r><w:rPr><w:i/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>Strona</w:t></w:r>
<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>_</w:t></w:r>
<w:r><w:rPr><w:i/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>name</w:t></w:r>
<w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>-hyperlink</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$newRange.InsertXML($xml)

# Step 4: InsertXML brought its own paragraph mark along with it, which
# pushed the previously-empty paragraph mark into a trailing empty
# paragraph. Remove that now-superfluous empty paragraph.
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanupRange = $d.Range($trailing.Range.Start - 1, $trailing.Range.End)
$cleanupRange.Delete()
